$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (subject ids) for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update CON row (row 2) for columns B:E
$ws.Range("B2").Value = 16.149999989999998
$ws.Range("C2").Value = 15.349999989999999
$ws.Range("D2").Value = 16.949999989999998
$ws.Range("E2").Value = 13.249999989999999

# Update STR row (row 3) for columns B:E
$ws.Range("B3").Value = 11.949999989999998
$ws.Range("C3").Value = 6.8499999899999997
$ws.Range("D3").Value = 19.749999989999999
$ws.Range("E3").Value = 23.79999999

# Update the selected range to reflect the new active selection
$ws.Range("B1:E3").Select() | Out-Null
